$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with refreshed crypto data
# D column forced to Text format to preserve exact numeric-string formatting
# (trailing zeros, decimal grouping, etc.) instead of being parsed as numbers.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.958.08"
$ws.Range("E2").Value = "  -1.83%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.742.16"
$ws.Range("E3").Value = "  -0.20%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.47%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.18"
$ws.Range("E5").Value = "  -4.83%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9998"
$ws.Range("E6").Value = "  -0.36%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5018"
$ws.Range("E7").Value = "  +5.72%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3577"
$ws.Range("E8").Value = "  +1.67%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "42.54"
$ws.Range("E9").Value = "  -0.30%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07269"
$ws.Range("E10").Value = "  -2.31%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.060"
$ws.Range("E11").Value = "  -1.08%  "
$ws.Range("E12").Value = "  -0.44%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.36"
$ws.Range("E13").Value = "  +0.53%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.998"
$ws.Range("E14").Value = "  -0.20%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.744.57"
$ws.Range("E15").Value = "  -0.05%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.869"
$ws.Range("E16").Value = "  -2.14%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "86.92"
$ws.Range("E17").Value = "  -5.21%  "
$ws.Range("E18").Value = "  -3.30%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06397"
$ws.Range("E19").Value = "  +0.47%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9998"
$ws.Range("E20").Value = "  -0.39%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.60"
$ws.Range("E21").Value = "  -0.31%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.738"
$ws.Range("E22").Value = "  -0.01%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "27.041.60"
$ws.Range("E23").Value = "  -1.82%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.22"
$ws.Range("E24").Value = "  +2.37%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.049"
$ws.Range("E25").Value = "  -4.62%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "154.73"
$ws.Range("E26").Value = "  -3.90%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.01"
$ws.Range("E27").Value = "  +1.02%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.942.59"
$ws.Range("E28").Value = "  -0.19%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.149"
$ws.Range("E29").Value = "  -1.48%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "121.45"
$ws.Range("E30").Value = "  +0.27%  "
$ws.Range("E31").Value = "  +0.82%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09531"
$ws.Range("E32").Value = "  +2.39%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.572"
$ws.Range("E33").Value = "  -1.37%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.394"
$ws.Range("E34").Value = "  -1.05%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02209"
$ws.Range("E35").Value = "  -1.35%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.05900"
$ws.Range("E36").Value = "  -0.44%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "11.12"
$ws.Range("E37").Value = "  -1.57%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.436"
$ws.Range("E38").Value = "  +0.08%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2001"
$ws.Range("E39").Value = "  -1.99%  "
$ws.Range("E40").Value = "  -1.18%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6037"
$ws.Range("E41").Value = "  -0.37%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9993"
$ws.Range("E42").Value = "  -0.40%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.119"
$ws.Range("E43").Value = "  -4.59%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "7.534"
$ws.Range("E44").Value = "  -2.11%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "12.86"
$ws.Range("E45").Value = "  -0.40%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.603"
$ws.Range("E46").Value = "  -3.24%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5660"
$ws.Range("E47").Value = "  -0.43%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "119.98"
$ws.Range("E48").Value = "  -2.06%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.863"
$ws.Range("E49").Value = "  -1.85%  "
$ws.Range("E50").Value = "  -1.81%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06670"
$ws.Range("E51").Value = "  -1.21%  "
